$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and date range) ---
$ws.Range("A8").Value = "Volume 30   Number  50"
$ws.Range("C9").Value = "Report Covering the Week  12/11/2023  Through  12/17/2023"

# --- Row 14-29 data table updates ---
# Row 14
$ws.Range("M14").Value = -46.153846153846
$ws.Range("N14").Value = -80
# Row 15
$ws.Range("G16").Copy($ws.Range("C15"))
$ws.Range("C15").Value = 2
$ws.Range("G16").Copy($ws.Range("D15"))
$ws.Range("D15").Value = 1
$ws.Range("H16").Copy($ws.Range("E15"))
$ws.Range("E15").Value = 100
$ws.Range("G16").Copy($ws.Range("F15"))
$ws.Range("F15").Value = 2
$ws.Range("H15").Value = 100
$ws.Range("I15").Value = 13
$ws.Range("J15").Value = 13
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 8.333333333333
$ws.Range("M15").Value = -23.529411764705
$ws.Range("N15").Value = -67.5
# Row 16
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 16
$ws.Range("H16").Value = 6.666666666666
$ws.Range("I16").Value = 271
$ws.Range("J16").Value = 231
$ws.Range("K16").Value = 17.316017316017
$ws.Range("L16").Value = 17.316017316017
$ws.Range("M16").Value = 22.072072072072
$ws.Range("N16").Value = -69.378531073446
# Row 17
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 10
$ws.Range("E17").Value = -30
$ws.Range("F17").Value = 27
$ws.Range("G17").Value = 36
$ws.Range("H17").Value = -25
$ws.Range("I17").Value = 370
$ws.Range("J17").Value = 433
$ws.Range("K17").Value = -14.549653579676
$ws.Range("L17").Value = 1.648351648351
$ws.Range("M17").Value = 65.178571428571
$ws.Range("N17").Value = -43.251533742331
# Row 18
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 50
$ws.Range("F18").Value = 8
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 114
$ws.Range("J18").Value = 180
$ws.Range("K18").Value = -36.666666666666
$ws.Range("L18").Value = 11.764705882352
$ws.Range("M18").Value = 7.547169811320
$ws.Range("N18").Value = -77.949709864603
# Row 19
$ws.Range("C19").Value = 4
$ws.Range("D19").Value = 4
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 25
$ws.Range("G19").Value = 42
$ws.Range("H19").Value = -40.476190476190
$ws.Range("I19").Value = 525
$ws.Range("J19").Value = 601
$ws.Range("K19").Value = -12.645590682196
$ws.Range("L19").Value = 36.363636363636
$ws.Range("M19").Value = 132.300884955752
$ws.Range("N19").Value = 38.157894736842
# Row 20
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 3
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = -50
$ws.Range("I20").Value = 85
$ws.Range("J20").Value = 85
$ws.Range("L20").Value = 39.344262295082
$ws.Range("M20").Value = 23.188405797101
$ws.Range("N20").Value = -74.164133738601
# Row 21
$ws.Range("C21").Value = 24
$ws.Range("D21").Value = 24
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 81
$ws.Range("G21").Value = 108
$ws.Range("H21").Value = -25
$ws.Range("I21").Value = 1385
$ws.Range("J21").Value = 1551
$ws.Range("K21").Value = -10.702772404900
$ws.Range("L21").Value = 18.274978650725
$ws.Range("M21").Value = 57.924743443557
$ws.Range("N21").Value = -51.198026779422
# Row 22
$ws.Range("D22").Value = 2
$ws.Range("C22").Copy($ws.Range("F22"))
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = -100
$ws.Range("J22").Value = 41
$ws.Range("K22").Value = -53.658536585365
$ws.Range("L22").Value = -42.424242424242
# Row 23
$ws.Range("C23").Value = 3
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 50
$ws.Range("G23").Value = 12
$ws.Range("H23").Value = -16.666666666666
$ws.Range("I23").Value = 156
$ws.Range("J23").Value = 180
$ws.Range("K23").Value = -13.333333333333
$ws.Range("L23").Value = -3.703703703703
$ws.Range("M23").Value = 56
# Row 24
$ws.Range("C24").Value = 15
$ws.Range("D24").Value = 26
$ws.Range("E24").Value = -42.307692307692
$ws.Range("F24").Value = 62
$ws.Range("G24").Value = 96
$ws.Range("H24").Value = -35.416666666666
$ws.Range("I24").Value = 1248
$ws.Range("J24").Value = 1323
$ws.Range("K24").Value = -5.668934240362
$ws.Range("L24").Value = 28.263103802672
$ws.Range("M24").Value = 30.954879328436
# Row 25
$ws.Range("C25").Value = 13
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = 44.444444444444
$ws.Range("F25").Value = 50
$ws.Range("G25").Value = 31
$ws.Range("H25").Value = 61.290322580645
$ws.Range("I25").Value = 532
$ws.Range("J25").Value = 540
$ws.Range("K25").Value = -1.481481481481
$ws.Range("L25").Value = 5.555555555555
$ws.Range("M25").Value = 10.373443983402
# Row 26
$ws.Range("G16").Copy($ws.Range("C26"))
$ws.Range("C26").Value = 2
$ws.Range("G16").Copy($ws.Range("D26"))
$ws.Range("D26").Value = 2
$ws.Range("H16").Copy($ws.Range("E26"))
$ws.Range("E26").Value = 0
$ws.Range("G16").Copy($ws.Range("F26"))
$ws.Range("F26").Value = 2
$ws.Range("G26").Value = 3
$ws.Range("H26").Value = -33.333333333333
$ws.Range("I26").Value = 23
$ws.Range("J26").Value = 29
$ws.Range("K26").Value = -20.689655172413
$ws.Range("L26").Value = -11.538461538461
# Row 27
$ws.Range("C27").Value = 4
$ws.Range("I27").Value = 75
$ws.Range("K27").Value = 20.967741935483
$ws.Range("L27").Value = 20.967741935483
# Row 28
$ws.Range("C22").Copy($ws.Range("C28"))
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = -100
$ws.Range("C22").Copy($ws.Range("F28"))
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = -100
$ws.Range("I28").Value = 29
$ws.Range("J28").Value = 43
$ws.Range("K28").Value = -32.558139534883
$ws.Range("L28").Value = -14.705882352941
$ws.Range("M28").Value = -27.5
$ws.Range("N28").Value = -62.337662337662
# Row 29
$ws.Range("C22").Copy($ws.Range("C29"))
$ws.Range("E29").Value = -100
$ws.Range("C22").Copy($ws.Range("F29"))
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = -100
$ws.Range("I29").Value = 27
$ws.Range("J29").Value = 31
$ws.Range("K29").Value = -12.903225806451
$ws.Range("L29").Value = 3.846153846153
$ws.Range("M29").Value = -15.625
$ws.Range("N29").Value = -62.5
